$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HW Q's")

$values = @(
  "*Given network connection diagram* Which device is the server? Which device is the client?",
  "Select all attributes that are provided in an IRIG-B-002 communication. (day, month, hour, minute, second - no year!)",
  "*Given engineering units scaling information* If the analog input measures 7.56mA, what will the output value be?",
  "*Given analog input resolution value and sensory input range* Evaluate the sensory resolution for the analog input.",
  "*Given motor's power in watts* Evaluate the horsepower rating of the specified motor.",
  "*Given pulse count resolution of rotary encoder* What is the smallest angular unit that this rotary encoder will record?",
  "What does MSB indicate? How about LSB?",
  "What level of current should motor contactors be rated for?",
  "What sort of cable (null-modem/straight-through) should be used for communications between an RTAC and relay?",
  "*Given voltage scaling information -PTR/scale/input-level* What should a technician expect to see as an output voltage?",
  "Suggest reasons why there are both AC and DC contact ratings for switches and relay contacts.",
  "*Given a link-budget diagram* What will the receive power be given an input power (some dB value)?"
)

$startRow = 18
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($startRow + $i, 4).Value = $values[$i]
}

$ws.Activate()
$ws.Range("D30").Select()
$excel.ActiveWindow.ScrollRow = 23
